$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns B:E to be treated as Text so numeric-looking / percent strings
# are not auto-converted to numbers by Excel when we assign them.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '41.179.60'
$ws.Range("E2").Value = '  -1.40%  '
$ws.Range("D3").Value = '2.173.81'
$ws.Range("E3").Value = '  -2.24%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '250.90'
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("E6").Value = '  -2.15%  '
$ws.Range("D7").Value = '66.03'
$ws.Range("E7").Value = '  -8.17%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = '0.575'
$ws.Range("E9").Value = '  -3.86%  '
$ws.Range("D10").Value = '58.99'
$ws.Range("E10").Value = '  +1.11%  '
$ws.Range("D11").Value = '36.18'
$ws.Range("E11").Value = '  -11.74%  '
$ws.Range("E12").Value = '  -3.64%  '
$ws.Range("E13").Value = '  -1.82%  '
$ws.Range("E14").Value = '  -5.37%  '
$ws.Range("D15").Value = '2.499.23'
$ws.Range("E15").Value = '  -2.16%  '
$ws.Range("E16").Value = '  -4.74%  '
$ws.Range("D17").Value = '0.841'
$ws.Range("E17").Value = '  -2.82%  '
$ws.Range("D18").Value = '2.176.93'
$ws.Range("E18").Value = '  -1.99%  '
$ws.Range("D19").Value = '41.091.23'
$ws.Range("E19").Value = '  -1.56%  '
$ws.Range("D20").Value = '0.0₃0943'
$ws.Range("E20").Value = '  -2.42%  '
$ws.Range("D21").Value = '71.50'
$ws.Range("E21").Value = '  -1.93%  '
$ws.Range("D22").Value = '6.04'
$ws.Range("E22").Value = '  -3.13%  '
$ws.Range("D23").Value = '229.87'
$ws.Range("E23").Value = '  -2.23%  '
$ws.Range("D24").Value = '2.00'
$ws.Range("E24").Value = '  -4.97%  '
$ws.Range("D25").Value = '3.80'
$ws.Range("E25").Value = '  -6.48%  '
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = '11.34'
$ws.Range("E27").Value = '  +5.77%  '
$ws.Range("D28").Value = '2.40'
$ws.Range("E28").Value = '  -5.42%  '
$ws.Range("D29").Value = '2.21'
$ws.Range("E29").Value = '  +0.88%  '
$ws.Range("D30").Value = '168.61'
$ws.Range("E30").Value = '  -1.56%  '
$ws.Range("D31").Value = '20.15'
$ws.Range("E31").Value = '  -2.99%  '
$ws.Range("D32").Value = '0.121'
$ws.Range("E32").Value = '  -2.91%  '
$ws.Range("D33").Value = '5.70'
$ws.Range("E33").Value = '  +1.57%  '
$ws.Range("D34").Value = '0.0747'
$ws.Range("E34").Value = '  +1.98%  '
$ws.Range("E35").Value = '  -3.53%  '
$ws.Range("D36").Value = '4.50'
$ws.Range("E36").Value = '  -4.81%  '
$ws.Range("E37").Value = '  -2.35%  '
$ws.Range("D38").Value = '24.38'
$ws.Range("E38").Value = '  -4.87%  '
$ws.Range("D39").Value = '0.0304'
$ws.Range("E39").Value = '  +0.96%  '
$ws.Range("D40").Value = '5.44'
$ws.Range("E40").Value = '  +11.58%  '
$ws.Range("E41").Value = '  -3.97%  '
$ws.Range("E42").Value = '  -8.08%  '
$ws.Range("B43").Value = 'MultiversX'
$ws.Range("C43").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D43").Value = '60.60'
$ws.Range("E43").Value = '  -9.13%  '
$ws.Range("B44").Value = 'Celestia'
$ws.Range("C44").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D44").Value = '11.30'
$ws.Range("E44").Value = '  -6.97%  '
$ws.Range("D45").Value = '8.46'
$ws.Range("E45").Value = '  -3.64%  '
$ws.Range("B46").Value = 'BinanceUSD'
$ws.Range("C46").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  -0.18%  '
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").Value = '0.188'
$ws.Range("E47").Value = '  -7.17%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.0988'
$ws.Range("E48").Value = '  -3.42%  '
$ws.Range("E49").Value = '  -2.84%  '
$ws.Range("B50").Value = 'SynthetixNetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D50").Value = '4.20'
$ws.Range("E50").Value = '  -9.59%  '
$ws.Range("B51").Value = 'TrustWalletToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D51").Value = '1.14'
$ws.Range("E51").Value = '  -4.63%  '

# Restore original (unformatted) cell style so the saved file matches the
# original workbook formatting (no explicit per-cell styles were present).
$ws.Range("B2:E51").ClearFormats()
